# This workbook contains a weekly price log for "Espinaca" (spinach) at
# "Vega Modelo de Temuco". The commit adds two new weekly records into the
# middle of the existing data table, which pushes all subsequent rows down
# (first by one row, then - after the second new record - by one more row).
#
# New record #1 is inserted at row 61 (everything that used to start at
# row 61 moves down to row 62).
# New record #2 is inserted at what is then row 102 (everything that used
# to start there moves down one more row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new weekly record (becomes row 61) ---
$ws.Rows.Item(61).Insert()

$ws.Range("A61").Value = 10
$ws.Range("B61").Value = "Vega Modelo de Temuco"
$ws.Range("C61").Value = "La Araucanía"
$ws.Range("D61").Value = 44664
$ws.Range("E61").Value = 9
$ws.Range("F61").Value = 100112012
$ws.Range("G61").Value = "Espinaca"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 40
$ws.Range("K61").Value = 9000
$ws.Range("L61").Value = 9000
$ws.Range("M61").Value = 9000
$ws.Range("N61").Value = "$/docena de atados"
$ws.Range("O61").Value = "Región de La Araucanía"
$ws.Range("P61").Value = 3000
$ws.Range("Q61").Value = 3
$ws.Range("R61").Value = "Hortaliza"

# --- Insert second new weekly record (becomes row 102) ---
$ws.Rows.Item(102).Insert()

$ws.Range("A102").Value = 10
$ws.Range("B102").Value = "Vega Modelo de Temuco"
$ws.Range("C102").Value = "La Araucanía"
$ws.Range("D102").Value = 44663
$ws.Range("E102").Value = 9
$ws.Range("F102").Value = 100112012
$ws.Range("G102").Value = "Espinaca"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 30
$ws.Range("K102").Value = 9000
$ws.Range("L102").Value = 9000
$ws.Range("M102").Value = 9000
$ws.Range("N102").Value = "$/docena de atados"
$ws.Range("O102").Value = "Región de La Araucanía"
$ws.Range("P102").Value = 3000
$ws.Range("Q102").Value = 3
$ws.Range("R102").Value = "Hortaliza"
